$newText = @'
questions = [
    {
        "title": "As a data analyst for a large grocery store, you have a dataset of customer purchases. Product IDs are 10 digits. The first three ID digits identify the product category, and the next three digits specify country and manufacturer. You want to use the product IDs to filter and aggregate the data by product. How would you change the Product ID field to allow for easy analysis of customer behavior based on geography and product category?",
        "ques_type": 2,
        "options": [
            "Convert to a text-based field such as string or character.",
            "Convert to a whole decimal number-based field such as integer.",
            "Convert to a field that supports decimal points.",
            "Keep it as-is and use any field that is already in the system."
        ],
        "score": "Convert to a text-based field such as string or character."
    },
    {
        "title": "You work as a data visualizer at a software company in the education space. Your team has put together the bar graph below to help busy school administrators see the total number of tutorials offered to their students every month and help them understand how the total is calculated.How should you modify the presentation of the data to achieve this?",
        "ques_type": 15,
        "options": [
            "Change the bar graph to a stacked bar graph and use colors to differentiate tutorial types.",
            "Add a brief description to each one of the categories shown.",
            "Change the graph to a line graph and chart over time",
            "Make the chart more colorful.",
            "Show all the values in a table over time."
        ],
        "score": [
            "Change the bar graph to a stacked bar graph and use colors to differentiate tutorial types.",
            "Add a brief description to each one of the categories shown."
        ]
    },
    {
        "title": "You are a product analyst in the consumer packaged goods industry working for one of the major manufacturers. Your manager has asked you to create a short write-up on customer behavior for a new product launch to be included in a brochure at a conference she is planning to attend. How should you start working on this project?",
        "ques_type": 2,
        "options": [
            "Collect data on customer demographics and past purchasing behavior.",
            "Start analyzing the product data and prepare a report intended for your ideal customers.",
            "Ask the manager who will be attending the conference what she aims to demonstrate in the write-up.",
            "Use industry standards to prepare a brief write-up for your manager."
        ],
        "score": "Ask the manager who will be attending the conference what she aims to demonstrate in the write-up."
    },
    {
        "title": "You are reviewing data for a national chain of coffee shops where prices, branding, and operations are set on a national level. You notice that the sales for a specific location close to a busy office park have decreased significantly compared to the same time last year. None of the other locations have experienced a similar decrease in sales. What hypothesis about the potential reason for this decrease in sales should you investigate first?",
        "ques_type": 2,
        "options": [
            "A decrease in foot traffic in the area.",
            "A change in consumer preferences towards tea.",
            "A decrease in the number of employees working at the location.",
            "An increase in the price of coffee at the location."
        ],
        "score": "A decrease in foot traffic in the area."
    }
]
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 2 (the shared-string cell lived there before)
$ws.Range("A2").EntireRow.Delete()

# Reset A1 formatting back to the default "Normal" style (removes bold font,
# thin border, and center/top alignment that the old header-like cell had)
$ws.Range("A1").Style = "Normal"

# Set A1 to the new (reformatted) text content
$ws.Range("A1").Value = $newText

# Restore the default (non-custom) row height after the long multi-line
# text would otherwise trigger an auto row-height override
$ws.Rows(1).AutoFit()
